# Update "想去人数" (interested-count) values for a handful of events.
# Sheet 1: 展览 (Exhibition)
# Sheet 2: 演出 (Performance)
# Sheet 4: 全部类型 (All Types) - mirrors the same events as sheets 1-3

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 109
$wsExhibit.Range("F3").Value = 417

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 26

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 109
$wsAll.Range("F4").Value = 417
$wsAll.Range("F8").Value = 26
